$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy styles from row 3 into row 4 first (so A4/G4 reuse the existing
# date / bool styles instead of Excel minting new ones), then overwrite
# with the new row's values.
$ws.Range("A3").Copy($ws.Range("A4"))
$ws.Range("G3").Copy($ws.Range("G4"))

$ws.Range("A4").Value = 42633.679097222222
$ws.Range("B4").Value = $false
$ws.Range("C4").Value = 9974
$ws.Range("D4").Value = 10000
$ws.Range("E4").Value = 19.32
$ws.Range("F4").Value = 19.22
$ws.Range("G4").Value = $false
$ws.Range("H4").Value = -0.52
$ws.Range("I4").Value = $false
